$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 98) taken from DGS's 2021/10/25 report.
# Column A holds the report date as text (sharing the same display format
# as the existing date cells above it), columns B:E hold the numeric data.
$dateFormat = $ws.Range("A97").NumberFormat

$ws.Range("A98").NumberFormat = "@"
$ws.Range("A98").Value = "2021/10/25"
$ws.Range("A98").NumberFormat = $dateFormat

$ws.Range("B98").Value = 92.4
$ws.Range("C98").Value = 92.8
$ws.Range("D98").Value = 1.06
$ws.Range("E98").Value = 1.06

# Move the selection to the next empty row, as in the saved workbook.
$ws.Range("A99").Select() | Out-Null
